$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.423.81'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '2.585.92'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.27'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.61'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("E8").Value = '  -7.40%  '
$ws.Range("D9").Value = '2.592.07'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.56'
$ws.Range("E10").Value = '  +6.49%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.346'
$ws.Range("E12").Value = '  +2.17%  '
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '3.039.00'
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("D15").Value = '60.434.19'
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.42'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("D18").Value = '2.589.72'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.02'
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.13'
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.01'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '0.0₃0843'
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.31'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.32'
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.57'
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.97'
$ws.Range("E35").Value = '  +2.13%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.861'
$ws.Range("E37").Value = '  +14.05%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("E41").Value = '  +2.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '296.03'
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0998'
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.614'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0556'
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.80'
$ws.Range("E47").Value = '  +3.55%  '
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("D51").Value = '1.999.78'
$ws.Range("E51").Value = '  +0.02%  '
